$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update prices in D29 and D30
$ws.Range("D29").Value = 520.458
$ws.Range("D30").Value = 353.073

# Re-apply merges (same ranges) in the order seen in the target workbook,
# which re-creates the merged-cell bookkeeping in the new order.
$mergeRanges = @(
    "A10:D10",
    "B38:C38",
    "A11:D11",
    "A12:D12",
    "B37:C37",
    "A1:D1",
    "B35:C35",
    "B39:C39",
    "B30:C30",
    "B28:C28",
    "A9:D9",
    "B36:C36",
    "B29:C29"
)
foreach ($r in $mergeRanges) {
    $ws.Range($r).MergeCells = $false
    $ws.Range($r).Merge()
}
